$d = $word.ActiveDocument

# The document has a paragraph whose entire content is a tracked
# insertion (w:ins) containing the text "insertion". The target edit
# removes that whole paragraph (including its paragraph mark), leaving
# "Test 1" immediately followed by "Test 2", while leaving the other
# pre-existing tracked change (the inserted paragraph mark on "Test 1")
# untouched.

# First, reject the tracked insertion of the run "insertion" so its
# text is actually removed rather than merely marked for deletion.
$d.Revisions(2).Reject()

# The paragraph mark of that now-empty paragraph was never itself a
# tracked insertion, so simply deleting it while TrackRevisions is on
# would just wrap it in a tracked deletion instead of removing it. Turn
# tracking off for this structural cleanup, delete the now-empty
# paragraph, then restore the TrackRevisions setting so the document's
# track-changes setting is preserved.
$d.TrackRevisions = $false
$d.Paragraphs(2).Range.Delete()
$d.TrackRevisions = $true
